$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Reposition the "H:M" quest/NPC/items content table (rows 14-20 -> 10-17)
#    Copy formats first (ascending source order is safe since destinations
#    are always <= source row and never revisited as a later source).
# ---------------------------------------------------------------------------
$hmPairs = @(
    @(14,10),
    @(15,11),
    @(16,13),
    @(17,14),
    @(18,15),
    @(19,16),
    @(20,17)
)
foreach ($pair in $hmPairs) {
    $src = $pair[0]
    $dst = $pair[1]
    $ws.Range("H$src`:M$src").Copy() | Out-Null
    $ws.Range("H$dst`:M$dst").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# New row 12 (F: Ruskie Shakel row) has no old counterpart - clone format from row 11
$ws.Range("H11:M11").Copy() | Out-Null
$ws.Range("H12:M12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Grow the "B:G" bordered box by one row: shift rows 13-18 -> 14-19.
#    Must process in DESCENDING source order to avoid clobbering a row
#    before it has been read.
# ---------------------------------------------------------------------------
$bgPairs = @(
    @(18,19),
    @(17,18),
    @(16,17),
    @(15,16),
    @(14,15),
    @(13,14)
)
foreach ($pair in $bgPairs) {
    $src = $pair[0]
    $dst = $pair[1]
    $ws.Range("B$src`:G$src").Copy() | Out-Null
    $ws.Range("B$dst`:G$dst").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Row 13 (C,D,E,F,G) must lose the box-border formatting it used to carry -
# only column B keeps its (unchanged) plain style there.
$ws.Range("C13:G13").ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 3) Move the "Key:" legend table from L22:M25 up to J5:L8 (gains a 3rd
#    "Fa = Farmable" column). Copy formats from the old legend cells.
# ---------------------------------------------------------------------------
$ws.Range("L22").Copy() | Out-Null
$ws.Range("J5").PasteSpecial(-4122) | Out-Null
$ws.Range("M22").Copy() | Out-Null
$ws.Range("K5:L5").PasteSpecial(-4122) | Out-Null

$ws.Range("L23").Copy() | Out-Null
$ws.Range("J6").PasteSpecial(-4122) | Out-Null
$ws.Range("M23").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("L24").Copy() | Out-Null
$ws.Range("L6").PasteSpecial(-4122) | Out-Null

$ws.Range("L23").Copy() | Out-Null
$ws.Range("J7").PasteSpecial(-4122) | Out-Null
$ws.Range("M24").Copy() | Out-Null
$ws.Range("K7:L7").PasteSpecial(-4122) | Out-Null

$ws.Range("L25").Copy() | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("M25").Copy() | Out-Null
$ws.Range("K8:L8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Clear now-unused leftover cells/old formatting from the pre-move layout.
# ---------------------------------------------------------------------------
$ws.Range("K14").Clear() | Out-Null
$ws.Range("H18:M20").Clear() | Out-Null
$ws.Range("L22:M25").Clear() | Out-Null
$ws.Range("B21:M25").ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 5) Write the actual cell values/formulas for the new layout.
# ---------------------------------------------------------------------------

# Key legend block (J5:L8)
$ws.Range("J5").Value = "Key:"
$ws.Range("J6").Value = "Q = Quest"
$ws.Range("K6").Value = "I = Item"
$ws.Range("L6").Value = "Fa = Farmable"
$ws.Range("J7").Value = "D = DROP"
$ws.Range("K7").Value = "G = GIVE"
$ws.Range("J8").Value = "E = Enemy"
$ws.Range("K8").Value = "F = Friendly"

# World/NPC/Items content table (H10:M17)
$ws.Range("H10").Value = "WORLD"
$ws.Range("I10").Value = "Room1"
$ws.Range("J10").Value = "Room2"
$ws.Range("K10").Value = "Room5"
$ws.Range("L10").Value = "Room3"
$ws.Range("M10").Value = "Room4"

$ws.Range("H11").Value = "NPC:"
$ws.Range("I11").Value = "F: Bel Drock"
$ws.Range("K11").Value = "F: Charlie"
$ws.Range("L11").Value = "E: FA: Altayere"
$ws.Range("M11").Value = "E: FA:(Diplodicus)"

$ws.Range("I12").Value = "F: Ruskie Shakel (After Q2)"

$ws.Range("H13").Value = "Items:"
$ws.Range("I13").Value = "Q1: Get QI1"
$ws.Range("K13").Value = "G:Sword"
$ws.Range("L13").Value = "DQI1: Shield"
$ws.Range("M13").Value = "D: QI3 Some Item"

$ws.Range("I14").Value = "Q3: Get XP"

$ws.Range("H17").Value = "LOCK:"
$ws.Range("J17").Value = "UNTIL Q1"
$ws.Range("M17").Value = "UNTIL Q3"

# New single-space cell
$ws.Range("F14").Value = " "

# New Q1/Q2/Q3 status block (H19:I21) - plain, unstyled cells
$ws.Range("H19").Value = "Q1"
$ws.Range("I19").Value = "Get Shield"
$ws.Range("H20").Value = "Q2"
$ws.Range("I20").Value = "Get to Lvl 3"
$ws.Range("H21").Value = "Q3"
$ws.Range("I21").Value = "Kill Diplodicus"

# ---------------------------------------------------------------------------
# 6) Merge cell moved from L22:M22 to J5:L5.
# ---------------------------------------------------------------------------
$ws.Range("L22:M22").UnMerge() | Out-Null
$ws.Range("J5:L5").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 7) Move the embedded Visio OLE object shape to its new anchor position.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 152.64159171998
$shp.Top = 123.093779527559
$shp.Width = 206.1103515625
$shp.Height = 174

# ---------------------------------------------------------------------------
# 8) Sheet view changes: zoom and selected cell.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("I13").Select() | Out-Null

Write-Output "done"
